# Round 3 data: fill in the "6.09-5.63" row (row 8) totals that were
# previously blank, and update the current selection to reflect where
# the user left off after entering the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 1546
$ws.Range("E8").Value = 369
$ws.Range("F8").Value = 3238

$ws.Range("F9").Select()
